$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "mid classes" (7:00-9:00 slots / duplicate lecture) were buggy entries
# that need to be removed from the routine:
#   Row 5  (TUE 7:00-9:00, 5CS024 Collaborative Development, Mr. Raj Shrestha)
#   Row 8  (WED 7:00-9:00, 5CS022 Human Computer Interaction, Mr. Ayush Shakya)
#   Row 12 (THU 9:30-11:30, 5CS020 Distributed and Cloud Systems Programming, Mr. Sumanta Silwal - Lecture)
#
# Deleting these rows shifts everything else up, producing the final
# 10-row table (header + 9 data rows) seen in the target workbook.
# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows.Item(12).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
